# Daily attendance processing - 2026-01-06 09:11:14
# Reorders the "Recorded By" (column G) comma-separated list of recorders
# by moving the last entry to the front (rotate right by one) for every
# row whose value contains more than one entry.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count + $usedRange.Row - 1

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)   # Column G = 7
    $value = $cell.Value2

    if ($value -ne $null -and $value -is [string] -and $value.Contains(",")) {
        $parts = $value.Split(",")
        for ($i = 0; $i -lt $parts.Length; $i++) {
            $parts[$i] = $parts[$i].Trim()
        }

        if ($parts.Length -gt 1) {
            $rotated = New-Object 'object[]' $parts.Length
            $rotated[0] = $parts[$parts.Length - 1]
            for ($i = 0; $i -lt $parts.Length - 1; $i++) {
                $rotated[$i + 1] = $parts[$i]
            }
            $newValue = [string]::Join(", ", $rotated)
            $cell.Value2 = $newValue
        }
    }
}
